$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Row 11: new entry
$ws.Cells.Item(10, 1).Copy()
$ws.Cells.Item(11, 1).PasteSpecial(-4122)
$ws.Cells.Item(11, 1).Value = 44154
$ws.Cells.Item(11, 2).Value = 4
$ws.Range("C11").Formula = "=C10+B11"
$ws.Cells.Item(11, 4).Value = "Implementierung bidirektionale Kommunikation zwischen Client und Server, Evaluation von libsigrok"

# Row 12: new entry
$ws.Cells.Item(10, 1).Copy()
$ws.Cells.Item(12, 1).PasteSpecial(-4122)
$ws.Cells.Item(12, 1).Value = 44156
$ws.Cells.Item(12, 2).Value = 1
$ws.Cells.Item(12, 4).Value = "Diskussion UI-Mockup"

# Update selection to I8 to mirror final author state
$ws.Range("I8").Select()

$wb.Save()
